# Update "想去人数" (F column) counts across sheets to match newly scraped data
# (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1794
$ws1.Range("F4").Value = 459
$ws1.Range("F7").Value = 633
$ws1.Range("F8").Value = 341
$ws1.Range("F9").Value = 1742
$ws1.Range("F10").Value = 369
$ws1.Range("F11").Value = 1427
$ws1.Range("F12").Value = 815
$ws1.Range("F13").Value = 340
$ws1.Range("F14").Value = 684
$ws1.Range("F15").Value = 12840
$ws1.Range("F16").Value = 12830
$ws1.Range("F17").Value = 960
$ws1.Range("F18").Value = 744
$ws1.Range("F22").Value = 572
$ws1.Range("F23").Value = 2013
$ws1.Range("F24").Value = 37
$ws1.Range("F25").Value = 10
$ws1.Range("F27").Value = 64
$ws1.Range("F28").Value = 254
$ws1.Range("F29").Value = 679

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 86
$ws2.Range("F6").Value = 18

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1794
$ws4.Range("F6").Value = 459
$ws4.Range("F11").Value = 633
$ws4.Range("F13").Value = 341
$ws4.Range("F14").Value = 1742
$ws4.Range("F15").Value = 369
$ws4.Range("F16").Value = 1427
$ws4.Range("F17").Value = 815
$ws4.Range("F18").Value = 340
$ws4.Range("F19").Value = 86
$ws4.Range("F20").Value = 684
$ws4.Range("F21").Value = 12840
$ws4.Range("F22").Value = 12830
$ws4.Range("F23").Value = 960
$ws4.Range("F24").Value = 744
$ws4.Range("F28").Value = 572
$ws4.Range("F29").Value = 18
$ws4.Range("F31").Value = 2013
$ws4.Range("F32").Value = 37
$ws4.Range("F33").Value = 10
$ws4.Range("F37").Value = 64
$ws4.Range("F38").Value = 254
$ws4.Range("F39").Value = 679
